$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "D1"
$ws.Range("B5").Value = 34
$ws.Range("C5").Value = 39.2000000000000028421709430404007434844970703125
$ws.Range("D5").Value = 37.2999999999999971578290569595992565155029296875
$ws.Range("E5").Value = 36.7999999999999971578290569595992565155029296875

$ws.Range("A6").Select()
